# Edit: split the "Criminal dresses"/" as a gay man" paragraph back into two
# runs inside a single paragraph (no longer interrupted by the _GoBack
# bookmark), add two new bullet paragraphs ("Criminal frequently visits gay
# bars" stays, "Suspect goes by an alias" is new), and relocate the
# _GoBack bookmark to its own empty paragraph immediately after the new
# "Suspect goes by an alias" line.

$d = $word.ActiveDocument

# --- Locate the start of the "Criminal dresses" paragraph -----------------
$startRng = $d.Content.Duplicate
$startFound = $startRng.Find.Execute("Criminal dresses", $true, $false, $false,
                                      $false, $false, $true, 1, $false, "", 0)
if (-not $startFound) {
    throw "Could not find 'Criminal dresses' in the document."
}
$startPos = $startRng.Start

# --- Locate the end of the "Criminal frequently visits gay bars" paragraph -
$endRng = $d.Content.Duplicate
$endFound = $endRng.Find.Execute("Criminal frequently visits gay bars", $true,
                                  $false, $false, $false, $false, $true, 1,
                                  $false, "", 0)
if (-not $endFound) {
    throw "Could not find 'Criminal frequently visits gay bars' in the document."
}
# Extend past the matched text to include the paragraph mark so the whole
# paragraph (not just its text) is part of the replaced range.
$endRng.MoveEnd(1, 1) | Out-Null
$endPos = $endRng.End

# --- Replace the two paragraphs with the new four-paragraph block ---------
$target = $d.Range($startPos, $endPos)

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="' + $wNs + '">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t>Criminal dresses</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> as a gay man</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:t>Criminal frequently visits gay bars</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:t>Suspect goes by an alias</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
              '<w:bookmarkEnd w:id="0"/>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$target.InsertXML($newXml)
